$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "Sat Feb 17 22:56:07 EST 2024"
$ws.Range("B4").Value = "Sat Feb 17 22:56:21 EST 2024"
